$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.145036666666667
$ws.Range("H2").Value = 3.43511
$ws.Range("I2").Value = 0.4953865629219574
$ws.Range("J2").Value = 0.4953865629219574
$ws.Range("M2").Value = 13.713764
$ws.Range("N2").Value = 41.141292
$ws.Range("O2").Value = 0.0901423721847377
$ws.Range("P2").Value = 0.0901423721847377
$ws.Range("Q2").Value = 15.70276261801333
$ws.Range("R2").Value = 141.32486356212
$ws.Range("S2").Value = 0.04465531993022907
$ws.Range("T2").Value = 0.04465531993022907

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.145036666666667
$ws.Range("H3").Value = 3.43511
$ws.Range("I3").Value = 0.4953865629219574
$ws.Range("J3").Value = 0.4953865629219574
$ws.Range("N3").Value = 84.55600199999999
$ws.Range("O3").Value = 0.1852659027513629
$ws.Range("P3").Value = 0.1852659027513629
$ws.Range("Q3").Value = 32.27324089224666
$ws.Range("R3").Value = 290.45916803022
$ws.Range("S3").Value = 0.09177823879063128
$ws.Range("T3").Value = 0.09177823879063128

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.145036666666667
$ws.Range("H4").Value = 3.43511
$ws.Range("I4").Value = 0.4953865629219574
$ws.Range("J4").Value = 0.4953865629219574
$ws.Range("M4").Value = 21.07704566666666
$ws.Range("N4").Value = 63.23113699999999
$ws.Range("O4").Value = 0.1385421898057586
$ws.Range("P4").Value = 0.1385421898057586
$ws.Range("Q4").Value = 24.13399011334111
$ws.Range("R4").Value = 217.2059110200699
$ws.Range("S4").Value = 0.06863193922755621
$ws.Range("T4").Value = 0.06863193922755621

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.145036666666667
$ws.Range("H5").Value = 3.43511
$ws.Range("I5").Value = 0.4953865629219574
$ws.Range("J5").Value = 0.4953865629219574
$ws.Range("M5").Value = 89.15834833333334
$ws.Range("N5").Value = 267.475045
$ws.Range("O5").Value = 0.5860495352581409
$ws.Range("P5").Value = 0.5860495352581407
$ws.Range("Q5").Value = 102.0895779811056
$ws.Range("R5").Value = 918.8062018299501
$ws.Range("S5").Value = 0.2903210649735409
$ws.Range("T5").Value = 0.2903210649735408

$ws.Range("G6").Value = 0.4713496666666666
$ws.Range("I6").Value = 0.2039238551060172
$ws.Range("J6").Value = 0.2039238551060172
$ws.Range("M6").Value = 13.713764
$ws.Range("N6").Value = 41.141292
$ws.Range("O6").Value = 0.0901423721847377
$ws.Range("P6").Value = 0.0901423721847377
$ws.Range("Q6").Value = 6.463978090145332
$ws.Range("R6").Value = 58.175802811308
$ws.Range("S6").Value = 0.01838218004431313
$ws.Range("T6").Value = 0.01838218004431313

$ws.Range("G7").Value = 0.4713496666666666
$ws.Range("I7").Value = 0.2039238551060172
$ws.Range("J7").Value = 0.2039238551060172
$ws.Range("N7").Value = 84.55600199999999
$ws.Range("O7").Value = 0.1852659027513629
$ws.Range("P7").Value = 0.1852659027513629
$ws.Range("Q7").Value = 13.28514778578866
$ws.Range("S7").Value = 0.03778013710875441
$ws.Range("T7").Value = 0.03778013710875441

$ws.Range("G8").Value = 0.4713496666666666
$ws.Range("I8").Value = 0.2039238551060172
$ws.Range("J8").Value = 0.2039238551060172
$ws.Range("M8").Value = 21.07704566666666
$ws.Range("N8").Value = 63.23113699999999
$ws.Range("O8").Value = 0.1385421898057586
$ws.Range("P8").Value = 0.1385421898057586
$ws.Range("Q8").Value = 9.934658449301441
$ws.Range("R8").Value = 89.41192604371298
$ws.Range("S8").Value = 0.02825205744001986
$ws.Range("T8").Value = 0.02825205744001986

$ws.Range("G9").Value = 0.4713496666666666
$ws.Range("I9").Value = 0.2039238551060172
$ws.Range("J9").Value = 0.2039238551060172
$ws.Range("M9").Value = 89.15834833333334
$ws.Range("N9").Value = 267.475045
$ws.Range("O9").Value = 0.5860495352581409
$ws.Range("P9").Value = 0.5860495352581407
$ws.Range("Q9").Value = 42.02475776746722
$ws.Range("R9").Value = 378.222819907205
$ws.Range("S9").Value = 0.1195094805129299
$ws.Range("T9").Value = 0.1195094805129298

$ws.Range("G10").Value = 0.6323219999999999
$ws.Range("H10").Value = 1.896966
$ws.Range("I10").Value = 0.2735666300991275
$ws.Range("J10").Value = 0.2735666300991275
$ws.Range("M10").Value = 13.713764
$ws.Range("N10").Value = 41.141292
$ws.Range("O10").Value = 0.0901423721847377
$ws.Range("P10").Value = 0.0901423721847377
$ws.Range("Q10").Value = 8.671514680007999
$ws.Range("R10").Value = 78.04363212007199
$ws.Range("S10").Value = 0.02465994498772001
$ws.Range("T10").Value = 0.02465994498772002

$ws.Range("G11").Value = 0.6323219999999999
$ws.Range("H11").Value = 1.896966
$ws.Range("I11").Value = 0.2735666300991275
$ws.Range("J11").Value = 0.2735666300991275
$ws.Range("N11").Value = 84.55600199999999
$ws.Range("O11").Value = 0.1852659027513629
$ws.Range("P11").Value = 0.1852659027513629
$ws.Range("Q11").Value = 17.822206765548
$ws.Range("R11").Value = 160.399860889932
$ws.Range("S11").Value = 0.05068256868796302
$ws.Range("T11").Value = 0.05068256868796302

$ws.Range("G12").Value = 0.6323219999999999
$ws.Range("H12").Value = 1.896966
$ws.Range("I12").Value = 0.2735666300991275
$ws.Range("J12").Value = 0.2735666300991275
$ws.Range("M12").Value = 21.07704566666666
$ws.Range("N12").Value = 63.23113699999999
$ws.Range("O12").Value = 0.1385421898057586
$ws.Range("P12").Value = 0.1385421898057586
$ws.Range("Q12").Value = 13.327479670038
$ws.Range("R12").Value = 119.947317030342
$ws.Range("S12").Value = 0.03790051999171507
$ws.Range("T12").Value = 0.03790051999171508

$ws.Range("G13").Value = 0.6323219999999999
$ws.Range("H13").Value = 1.896966
$ws.Range("I13").Value = 0.2735666300991275
$ws.Range("J13").Value = 0.2735666300991275
$ws.Range("M13").Value = 89.15834833333334
$ws.Range("N13").Value = 267.475045
$ws.Range("O13").Value = 0.5860495352581409
$ws.Range("P13").Value = 0.5860495352581407
$ws.Range("Q13").Value = 56.37678513482999
$ws.Range("R13").Value = 507.39106621347
$ws.Range("S13").Value = 0.1603235964317294
$ws.Range("T13").Value = 0.1603235964317294

$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.062692
$ws.Range("H14").Value = 0.188076
$ws.Range("I14").Value = 0.02712295187289783
$ws.Range("J14").Value = 0.02712295187289783
$ws.Range("M14").Value = 13.713764
$ws.Range("N14").Value = 41.141292
$ws.Range("O14").Value = 0.0901423721847377
$ws.Range("P14").Value = 0.0901423721847377
$ws.Range("Q14").Value = 0.859743292688
$ws.Range("R14").Value = 7.737689634192
$ws.Range("S14").Value = 0.002444927222475485
$ws.Range("T14").Value = 0.002444927222475485

$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.062692
$ws.Range("H15").Value = 0.188076
$ws.Range("I15").Value = 0.02712295187289783
$ws.Range("J15").Value = 0.02712295187289783
$ws.Range("N15").Value = 84.55600199999999
$ws.Range("O15").Value = 0.1852659027513629
$ws.Range("P15").Value = 0.1852659027513629
$ws.Range("Q15").Value = 1.766994959128
$ws.Range("R15").Value = 15.902954632152
$ws.Range("S15").Value = 0.005024958164014185
$ws.Range("T15").Value = 0.005024958164014185

$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.062692
$ws.Range("H16").Value = 0.188076
$ws.Range("I16").Value = 0.02712295187289783
$ws.Range("J16").Value = 0.02712295187289783
$ws.Range("M16").Value = 21.07704566666666
$ws.Range("N16").Value = 63.23113699999999
$ws.Range("O16").Value = 0.1385421898057586
$ws.Range("P16").Value = 0.1385421898057586
$ws.Range("Q16").Value = 1.321362146934666
$ws.Range("R16").Value = 11.892259322412
$ws.Range("S16").Value = 0.003757673146467467
$ws.Range("T16").Value = 0.003757673146467467

$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.062692
$ws.Range("H17").Value = 0.188076
$ws.Range("I17").Value = 0.02712295187289783
$ws.Range("J17").Value = 0.02712295187289783
$ws.Range("M17").Value = 89.15834833333334
$ws.Range("N17").Value = 267.475045
$ws.Range("O17").Value = 0.5860495352581409
$ws.Range("P17").Value = 0.5860495352581407
$ws.Range("Q17").Value = 5.589515173713333
$ws.Range("R17").Value = 50.30563656342
$ws.Range("S17").Value = 0.01589539333994069
$ws.Range("T17").Value = 0.01589539333994069
